{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst oldText = \"Pas bcp de doc, difficile d'acceder au propri\u00e9t\u00e9 et a leur value (celle qui ont \u00e9t\u00e9 modifi\u00e9)\";\nconst newText = \"Pas beaucoup de doc, difficile d'acc\u00e9der aux propri\u00e9t\u00e9s et \u00e0 leur valeur (celle qui ont \u00e9t\u00e9 modifi\u00e9)\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === oldText) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\ntarget.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Pas bcp de doc, difficile d'acceder au propri\u00e9t\u00e9 et a leur value (celle qui ont \u00e9t\u00e9 modifi\u00e9)\"\n$newText = \"Pas beaucoup de doc, difficile d'acc\u00e9der aux propri\u00e9t\u00e9s et \u00e0 leur valeur (celle qui ont \u00e9t\u00e9 modifi\u00e9)\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $oldText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $start = $target.Range.Start\n    $end = $start + $oldText.Length\n    $repRange = $d.Range($start, $end)\n    $repRange.Text = $newText\n}\n"}
